$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the formatting from an existing
# header cell (H1) so they share the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I (I0) and J (IF), rows 2-31
$data = @(
    @(5, 7),
    @(6, 7),
    @(8, 8),
    @(7, 8),
    @(8, 9),
    @(3, 4),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(6, 6),
    @(4, 7),
    @(9, 9),
    @(5, 5),
    @(1, 2),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(6, 7),
    @(12, 12),
    @(9, 10),
    @(5, 7),
    @(8, 9),
    @(5, 6),
    @(3, 5),
    @(1, 3),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(7, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
